# Weekly update: insert a new price-report row for "Haba" (Terminal
# Hortofrutícola Agro Chillán) above the existing history, pushing all
# prior rows down by one. This mirrors how the source feed prepends the
# newest week's observation to the top of the data block (row 21, right
# below the still-unmoved rows 2-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 21:40 down to 22:41, opening up a blank row 21.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the latest weekly observation.
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44789
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112026
$ws.Range("G21").Value = "Haba"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 60
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 11500
$ws.Range("N21").Value = "`$/saco 25 kilos"
$ws.Range("O21").Value = "Provincia del Elquí"
$ws.Range("P21").Value = 460
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
